$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# Helper: build a minimal single-part WordprocessingML package used with
# Range.InsertXML to set the exact contents of a table-cell paragraph.
# ---------------------------------------------------------------------------
function New-CellXml([string]$innerParagraphsXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerParagraphsXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark from the end of the "Get all course
#    details..." question paragraph (it moves further down the document).
# ---------------------------------------------------------------------------
$qCell = $t.Rows.Item(46).Cells.Item(1)

$noBookmarkXml = New-CellXml @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr>
<w:ind w:left="454"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Get all course details which had started on ‘2016-02-01’.</w:t></w:r>
</w:p>
'@

$qCell.Range.InsertXML($noBookmarkXml)

# ---------------------------------------------------------------------------
# 2) Append six new rows at the end of the table.
#
#    Rows.Add() (no args) always appends after the table's current LAST row
#    and clones that last row's trPr/tcPr. Rows.Add(beforeRow) inserts just
#    before "beforeRow" and clones ITS trPr/tcPr instead. To avoid the two
#    "tall, vertically centered" answer rows (C and F) leaking their
#    w:trHeight / w:vAlign onto a neighbouring plain row, first append all
#    of the plain rows (A, B, D, E) while the table's tail is still plain,
#    then splice C in before D (cloning D's plain formatting) and append F
#    last (cloning E's plain formatting) - fixing up height/alignment only
#    by ADDING the properties, never by trying to remove them again.
#
#    Row-reference variables can go stale once later Add() calls reshuffle
#    the table, so every row is re-fetched by its final absolute index
#    before being touched again.
# ---------------------------------------------------------------------------

$base = $t.Rows.Count

$null = $t.Rows.Add()                 # becomes row A
$null = $t.Rows.Add()                 # becomes row B
$rowDHandle = $t.Rows.Add()           # becomes row D
$null = $t.Rows.Add()                 # becomes row E
$null = $t.Rows.Add($rowDHandle)      # becomes row C (spliced before D)
$null = $t.Rows.Add()                 # becomes row F

$idxA = $base + 1
$idxB = $base + 2
$idxC = $base + 3
$idxD = $base + 4
$idxE = $base + 5
$idxF = $base + 6

# --- Row A: blank spacer row -------------------------------------------------
$xmlA = New-CellXml @'
<w:p>
<w:pPr>
<w:ind w:left="94"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
</w:pPr>
</w:p>
'@
$t.Rows.Item($idxA).Cells.Item(1).Range.InsertXML($xmlA)

# --- Row B: question - "Get module names..." --------------------------------
$xmlB = New-CellXml @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr>
<w:ind w:left="454"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Get module names which are taught in ‘PG-DAC’ course.</w:t></w:r>
</w:p>
'@
$t.Rows.Item($idxB).Cells.Item(1).Range.InsertXML($xmlB)

# --- Row C: answer SQL for Row B --------------------------------------------
$rowC = $t.Rows.Item($idxC)
$rowC.Height = 27.5
$rowC.Cells.Item(1).VerticalAlignment = 1
$xmlC = New-CellXml @'
<w:p>
<w:pPr>
<w:ind w:left="454"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="5B9BD5" w:themeColor="accent1"/><w:sz w:val="20"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="5B9BD5" w:themeColor="accent1"/><w:sz w:val="20"/></w:rPr><w:t>select course.name, modules.name from course, modules, course_modules where course.id = course_modules.courseID and modules.ID = course_modules.moduleID and course.name='PG-DAC';</w:t></w:r>
</w:p>
'@
$t.Rows.Item($idxC).Cells.Item(1).Range.InsertXML($xmlC)

# --- Row D: blank spacer row -------------------------------------------------
$xmlD = New-CellXml @'
<w:p>
<w:pPr>
<w:ind w:left="94"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
</w:pPr>
</w:p>
'@
$t.Rows.Item($idxD).Cells.Item(1).Range.InsertXML($xmlD)

# --- Row E: question - "Display how many modules..." (with _GoBack bookmark)
$xmlE = New-CellXml @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr>
<w:ind w:left="454"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Display how many modules are taught </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">in </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>each course.</w:t></w:r>
</w:p>
'@
$t.Rows.Item($idxE).Cells.Item(1).Range.InsertXML($xmlE)

# --- Row F: answer SQL for Row E --------------------------------------------
$rowF = $t.Rows.Item($idxF)
$rowF.Height = 27.5
$rowF.Cells.Item(1).VerticalAlignment = 1
$xmlF = New-CellXml @'
<w:p>
<w:pPr>
<w:ind w:left="454"/>
<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="5B9BD5" w:themeColor="accent1"/><w:sz w:val="20"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:color w:val="5B9BD5" w:themeColor="accent1"/><w:sz w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>select course.name, count(modules.name) from course, modules, course_modules where course.id = course_modules.courseID and modules.ID = course_modules.moduleID group by course.name;</w:t></w:r>
</w:p>
'@
$t.Rows.Item($idxF).Cells.Item(1).Range.InsertXML($xmlF)

Write-Output ("Final row count: " + $t.Rows.Count)
